# Update values in the KNN imputed result data.
# Commit message: "Update Name of Algo" (values recomputed for the algorithm run)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.656000000000001
$ws.Range("B3").Value = 6.701000000000001
$ws.Range("D5").Value = -8.293000000000001
$ws.Range("E5").Value = 12.94
$ws.Range("E9").Value = 13.152
$ws.Range("E11").Value = 13.132
$ws.Range("B14").Value = 6.449
$ws.Range("B21").Value = 6.6
$ws.Range("E21").Value = 12.694
$ws.Range("B23").Value = 6.610000000000001
$ws.Range("B25").Value = 6.159000000000001
